$wb = $excel.ActiveWorkbook

# Scheduled runner refresh: updated market price snapshots (H/I/J/K/L columns)
# and recomputed profit deltas (M/N columns) across all 8 sheets.


$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 88.625
$ws.Range("I8").Value = 88.625
$ws.Range("K8").Value = 265.875
$ws.Range("M8").Value = -126.875
# Row 98
$ws.Range("H98").Value = 2609.6553
$ws.Range("I98").Value = 2486.7917
$ws.Range("K98").Value = 2486.7917
$ws.Range("M98").Value = -988.7917000000002
# Row 112
$ws.Range("H112").Value = 7881.4443
$ws.Range("J112").Value = 8169.4707
$ws.Range("L112").Value = 24508.4121
$ws.Range("N112").Value = -26724.4121
# Row 122
$ws.Range("H122").Value = 2609.6553
$ws.Range("I122").Value = 2486.7917
$ws.Range("K122").Value = 7460.375100000001
$ws.Range("M122").Value = -5010.375100000001
# Row 132
$ws.Range("H132").Value = 34486720
$ws.Range("I132").Value = 34486720
$ws.Range("K132").Value = 103460160
$ws.Range("M132").Value = -103457630
# Row 141
$ws.Range("H141").Value = 20118.125
$ws.Range("J141").Value = 41910
$ws.Range("L141").Value = 125730
$ws.Range("N141").Value = -136090

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5487.868
$ws.Range("I32").Value = 3950.442
$ws.Range("K32").Value = 3950.442
$ws.Range("M32").Value = -3663.442
# Row 122
$ws.Range("H122").Value = 1044377.56
$ws.Range("I122").Value = 1283.6666
$ws.Range("J122").Value = 1897818
$ws.Range("K122").Value = 3850.9998
$ws.Range("L122").Value = 5693454
$ws.Range("M122").Value = -1400.9998
$ws.Range("N122").Value = -5698354
# Row 132
$ws.Range("H132").Value = 4353.8
$ws.Range("I132").Value = 7059.579
$ws.Range("J132").Value = 2376.5
$ws.Range("K132").Value = 21178.737
$ws.Range("L132").Value = 7129.5
$ws.Range("M132").Value = -18648.737
$ws.Range("N132").Value = -12189.5

$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 200
$ws.Range("I11").Value = 200
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 200
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -60
$ws.Range("N11").ClearContents()
# Row 12
$ws.Range("H12").Value = 82
$ws.Range("J12").Value = 99
$ws.Range("L12").Value = 99
$ws.Range("N12").Value = -435
# Row 22
$ws.Range("H22").Value = 888.92
$ws.Range("I22").Value = 843.8095
$ws.Range("J22").Value = 1125.75
$ws.Range("K22").Value = 843.8095
$ws.Range("L22").Value = 1125.75
$ws.Range("M22").Value = -670.8095
$ws.Range("N22").Value = -1471.75
# Row 25
$ws.Range("H25").Value = 796
$ws.Range("I25").Value = 796
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 796
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -561
$ws.Range("N25").ClearContents()
# Row 37
$ws.Range("H37").Value = 4983.2
$ws.Range("I37").Value = 1463
$ws.Range("J37").Value = 7330
$ws.Range("K37").Value = 1463
$ws.Range("L37").Value = 7330
$ws.Range("M37").Value = -1326
$ws.Range("N37").Value = -7604
# Row 94
$ws.Range("H94").Value = 9092720
$ws.Range("I94").Value = 18183138
$ws.Range("J94").Value = 2302.2
$ws.Range("K94").Value = 18183138
$ws.Range("L94").Value = 2302.2
$ws.Range("M94").Value = -18182687
$ws.Range("N94").Value = -3204.2
# Row 99
$ws.Range("H99").Value = 4467174.5
$ws.Range("I99").Value = 6496216
$ws.Range("K99").Value = 6496216
$ws.Range("M99").Value = -6494718
# Row 134
$ws.Range("H134").Value = 14314.875
$ws.Range("I134").Value = 14398.789
$ws.Range("K134").Value = 43196.367
$ws.Range("M134").Value = -40661.367
# Row 135
$ws.Range("H135").Value = 45000
$ws.Range("J135").Value = 45000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -55140
# Row 140
$ws.Range("H140").Value = 83874.75
$ws.Range("J140").Value = 83874.75
$ws.Range("L140").Value = 83874.75
$ws.Range("N140").Value = -94234.75

$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
# Row 58
$ws.Range("H58").Value = 3524.5293
$ws.Range("I58").Value = 2890
$ws.Range("K58").Value = 2890
$ws.Range("M58").Value = -2687
# Row 134
$ws.Range("H134").Value = 9069.773999999999
$ws.Range("I134").Value = 6884
$ws.Range("J134").Value = 13043.909
$ws.Range("K134").Value = 20652
$ws.Range("L134").Value = 39131.727
$ws.Range("M134").Value = -18117
$ws.Range("N134").Value = -44201.727
# Row 136
$ws.Range("H136").Value = 3524.5293
$ws.Range("I136").Value = 2890
$ws.Range("K136").Value = 8670
$ws.Range("M136").Value = -6120
# Row 140
$ws.Range("H140").Value = 79000
$ws.Range("J140").Value = 79000
$ws.Range("L140").Value = 79000
$ws.Range("N140").Value = -89360

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 49467.777
$ws.Range("J12").Value = 111
$ws.Range("L12").Value = 333
$ws.Range("N12").Value = -679
# Row 37
$ws.Range("H37").Value = 68157
$ws.Range("J37").Value = 68157
$ws.Range("L37").Value = 204471
$ws.Range("N37").Value = -204695
# Row 120
$ws.Range("H120").Value = 15928.571
$ws.Range("J120").Value = 20300
$ws.Range("L120").Value = 60900
$ws.Range("N120").Value = -70576
# Row 122
$ws.Range("H122").Value = 600.6
$ws.Range("J122").Value = 226
$ws.Range("L122").Value = 2034
$ws.Range("N122").Value = -6934
# Row 128
$ws.Range("H128").Value = 179874.75
$ws.Range("I128").Value = 179874.75
$ws.Range("K128").Value = 539624.25
$ws.Range("M128").Value = -534644.25
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 28
$ws.Range("H28").Value = 19999.875
$ws.Range("I28").Value = 19999
$ws.Range("K28").Value = 19999
$ws.Range("M28").Value = -19807
# Row 102
$ws.Range("H102").Value = 3022997.2
$ws.Range("I102").Value = 4116642.2
$ws.Range("K102").Value = 4116642.2
$ws.Range("M102").Value = -4115020.2
# Row 107
$ws.Range("H107").Value = 1271.2858
$ws.Range("I107").Value = 1520
$ws.Range("J107").Value = 649.5
$ws.Range("K107").Value = 1520
$ws.Range("L107").Value = 649.5
$ws.Range("M107").Value = 400
$ws.Range("N107").Value = -4489.5
# Row 133
$ws.Range("H133").Value = 109997
$ws.Range("J133").Value = 109997
$ws.Range("L133").Value = 109997
$ws.Range("N133").Value = -120117

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 40657.48
$ws.Range("I22").Value = 66869.266
$ws.Range("J22").Value = 1339.8
$ws.Range("K22").Value = 66869.266
$ws.Range("L22").Value = 1339.8
$ws.Range("M22").Value = -66574.266
$ws.Range("N22").Value = -1929.8
# Row 27
$ws.Range("H27").Value = 40657.48
$ws.Range("I27").Value = 66869.266
$ws.Range("J27").Value = 1339.8
$ws.Range("K27").Value = 66869.266
$ws.Range("L27").Value = 1339.8
$ws.Range("M27").Value = -66762.266
$ws.Range("N27").Value = -1553.8
# Row 46
$ws.Range("H46").Value = 4584
$ws.Range("J46").Value = 10666.333
$ws.Range("L46").Value = 10666.333
$ws.Range("N46").Value = -11042.333
# Row 55
$ws.Range("H55").Value = 1480.0303
$ws.Range("I55").Value = 1582.4667
$ws.Range("J55").Value = 1394.6666
$ws.Range("K55").Value = 1582.4667
$ws.Range("L55").Value = 1394.6666
$ws.Range("M55").Value = -1409.4667
$ws.Range("N55").Value = -1740.6666
# Row 74
$ws.Range("H74").Value = 4197
$ws.Range("I74").Value = 4197
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4197
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3199
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 4197
$ws.Range("I77").Value = 4197
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 12591
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -7599
$ws.Range("N77").ClearContents()
# Row 82
$ws.Range("H82").Value = 2316334.5
$ws.Range("I82").Value = 3969855.5
$ws.Range("J82").Value = 1404.8
$ws.Range("K82").Value = 3969855.5
$ws.Range("L82").Value = 1404.8
$ws.Range("M82").Value = -3969494.5
$ws.Range("N82").Value = -2126.8
# Row 85
$ws.Range("H85").Value = 2316334.5
$ws.Range("I85").Value = 3969855.5
$ws.Range("J85").Value = 1404.8
$ws.Range("K85").Value = 3969855.5
$ws.Range("L85").Value = 1404.8
$ws.Range("M85").Value = -3968607.5
$ws.Range("N85").Value = -3900.8
# Row 97
$ws.Range("H97").Value = 2335.5
$ws.Range("J97").Value = 2335.5
$ws.Range("L97").Value = 2335.5
$ws.Range("N97").Value = -4317.5
# Row 105
$ws.Range("H105").Value = 44400
$ws.Range("J105").Value = 44400
$ws.Range("L105").Value = 44400
$ws.Range("N105").Value = -51388
# Row 132
$ws.Range("H132").Value = 11736.263
$ws.Range("I132").Value = 12599.267
$ws.Range("K132").Value = 37797.801
$ws.Range("M132").Value = -35267.801

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 8690.579
$ws.Range("I136").Value = 8796.467000000001
$ws.Range("K136").Value = 26389.401
$ws.Range("M136").Value = -23839.401
